$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Git Basic Commands")
$table = $ws.ListObjects.Item("Table1")

# Add two new rows of Git commands to the cheatsheet (new table rows grow
# the table range and the sheet automatically)
$newRow1 = $table.ListRows.Add()
$newRow1.Range.Cells(1, 1).Value = "git clone {repository_name_on_github}"
$newRow1.Range.Cells(1, 2).Value = "Clones a repository from GitHub and moves it to your local machine. Once you've done this, you can then put it on your own GitHub repository and mess around with it there"

$newRow2 = $table.ListRows.Add()
$newRow2.Range.Cells(1, 1).Value = "git branch -D {branch_name}"
$newRow2.Range.Cells(1, 2).Value = "Delete one of the git branches. Important to use a capital D"

# Row 14 needs a taller row height to fit the wrapped text (matches other
# multi-line rows in the sheet, e.g. row 8/9/13 which also use ht=45)
$ws.Rows.Item(14).RowHeight = 45

# Update the selection to match the target state
$ws.Range("B9").Select()
